$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-06-10T16:07:54"
$ws.Range("U4").Value = 46.18
$ws.Range("V4").Value = 34.05
$ws.Range("X4").Value = 34.14
$ws.Range("Y4").Value = 29.55
$ws.Range("Z4").Value = 14.93
$ws.Range("U6").Value = -1.06
$ws.Range("W6").Value = -0.5
$ws.Range("Y6").Value = -0.12
$ws.Range("Z6").Value = 0
$ws.Range("W8").Value = -131.42
$ws.Range("X8").Value = 4.21
$ws.Range("T9").Value = 46.73
$ws.Range("U9").Value = 47.24
$ws.Range("V9").Value = 35.09
$ws.Range("Y9").Value = 30.97
$ws.Range("Z9").Value = 15.61
$ws.Range("T11").Value = -0.33
$ws.Range("U11").Value = 0
$ws.Range("V11").Value = 0.46
$ws.Range("W11").Value = 0.52
$ws.Range("X11").Value = 1.38
$ws.Range("Y11").Value = 1.3
$ws.Range("Z11").Value = 0.6899999999999999
$ws.Range("W13").Value = -132.44
$ws.Range("X13").Value = 4.21
$ws.Range("T14").Value = 46.73
$ws.Range("U14").Value = 47.24
$ws.Range("V14").Value = 35.09
$ws.Range("X14").Value = 35.58
$ws.Range("Y14").Value = 31
$ws.Range("Z14").Value = 15.61
$ws.Range("T16").Value = -0.33
$ws.Range("U16").Value = 0
$ws.Range("V16").Value = 0.46
$ws.Range("W16").Value = 0.52
$ws.Range("X16").Value = 1.38
$ws.Range("Y16").Value = 1.33
$ws.Range("Z16").Value = 0.6899999999999999
$ws.Range("W18").Value = -132.44
$ws.Range("X18").Value = 4.21
$ws.Range("U19").Value = 46.5
$ws.Range("V19").Value = 34.42
$ws.Range("Y19").Value = 29.91
$ws.Range("Z19").Value = 15.12
$ws.Range("U21").Value = -0.74
$ws.Range("W21").Value = -0.16
$ws.Range("X21").Value = 0.33
$ws.Range("Y21").Value = 0.24
$ws.Range("Z21").Value = 0.2
$ws.Range("W23").Value = -131.76
$ws.Range("X23").Value = 4.21
$ws.Range("U24").Value = 46.5
$ws.Range("V24").Value = 34.42
$ws.Range("Y24").Value = 29.91
$ws.Range("Z24").Value = 15.12
$ws.Range("U26").Value = -0.74
$ws.Range("W26").Value = -0.16
$ws.Range("X26").Value = 0.33
$ws.Range("Y26").Value = 0.24
$ws.Range("Z26").Value = 0.2
$ws.Range("W28").Value = -131.76
$ws.Range("X28").Value = 4.21
$ws.Range("U29").Value = 46.92
$ws.Range("V29").Value = 34.8
$ws.Range("Y29").Value = 30.27
$ws.Range("Z29").Value = 15.33
$ws.Range("U31").Value = -0.33
$ws.Range("V31").Value = 0.17
$ws.Range("W31").Value = 0.23
$ws.Range("X31").Value = 0.8
$ws.Range("Y31").Value = 0.61
$ws.Range("Z31").Value = 0.4
$ws.Range("W33").Value = -132.14
$ws.Range("X33").Value = 4.21
$ws.Range("T34").Value = 46.73
$ws.Range("U34").Value = 47.48
$ws.Range("V34").Value = 35.3
$ws.Range("X34").Value = 35.91
$ws.Range("Y34").Value = 31.4
$ws.Range("Z34").Value = 15.75
$ws.Range("T36").Value = -0.33
$ws.Range("U36").Value = 0.24
$ws.Range("V36").Value = 0.67
$ws.Range("W36").Value = 0.75
$ws.Range("X36").Value = 1.71
$ws.Range("Y36").Value = 1.73
$ws.Range("Z36").Value = 0.82
$ws.Range("W38").Value = -132.67
$ws.Range("X38").Value = 4.21
$ws.Range("U39").Value = 46.18
$ws.Range("V39").Value = 34.05
$ws.Range("X39").Value = 34.14
$ws.Range("Y39").Value = 29.55
$ws.Range("Z39").Value = 14.93
$ws.Range("U41").Value = -1.06
$ws.Range("W41").Value = -0.5
$ws.Range("Y41").Value = -0.12
$ws.Range("Z41").Value = 0
$ws.Range("W43").Value = -131.42
$ws.Range("X43").Value = 4.21
$ws.Range("T44").Value = 46.45
$ws.Range("U44").Value = 46.64
$ws.Range("V44").Value = 34.02
$ws.Range("X44").Value = 33.99
$ws.Range("Y44").Value = 29.32
$ws.Range("Z44").Value = 14.79
$ws.Range("T46").Value = -0.6
$ws.Range("U46").Value = -0.61
$ws.Range("W46").Value = -0.41
$ws.Range("Y46").Value = -0.35
$ws.Range("Z46").Value = -0.13
$ws.Range("W48").Value = -131.51
$ws.Range("X48").Value = 4.21
$ws.Range("U49").Value = 43.66
$ws.Range("V49").Value = 31.86
$ws.Range("X49").Value = 31.62
$ws.Range("Y49").Value = 29.17
$ws.Range("Z49").Value = 14.63
$ws.Range("U51").Value = -3.58
$ws.Range("W51").Value = -2.39
$ws.Range("Y51").Value = -0.5
$ws.Range("Z51").Value = -0.29
$ws.Range("W53").Value = -129.53
$ws.Range("X53").Value = 4.21
$ws.Range("T54").Value = 43.9
$ws.Range("U54").Value = 44.24
$ws.Range("V54").Value = 32.45
$ws.Range("X54").Value = 32.96
$ws.Range("Y54").Value = 29.43
$ws.Range("Z54").Value = 14.94
$ws.Range("T56").Value = -3.16
$ws.Range("U56").Value = -3.01
$ws.Range("V56").Value = -2.17
$ws.Range("W56").Value = -1.49
$ws.Range("Y56").Value = -0.24
$ws.Range("Z56").Value = 0.01
$ws.Range("W58").Value = -130.43
$ws.Range("X58").Value = 4.21
$ws.Range("T59").Value = 47.97
$ws.Range("U59").Value = 48.06
$ws.Range("V59").Value = 34.87
$ws.Range("X59").Value = 34.84
$ws.Range("Y59").Value = 30.12
$ws.Range("Z59").Value = 15.12
$ws.Range("T61").Value = 0.91
$ws.Range("U61").Value = 0.82
$ws.Range("W61").Value = 0.52
$ws.Range("Y61").Value = 0.45
$ws.Range("Z61").Value = 0.2
$ws.Range("W63").Value = -132.44
$ws.Range("X63").Value = 4.21
$ws.Range("T64").Value = 48.71
$ws.Range("U64").Value = 48.75
$ws.Range("V64").Value = 35.41
$ws.Range("X64").Value = 35.29
$ws.Range("Y64").Value = 30.52
$ws.Range("Z64").Value = 15.29
$ws.Range("T66").Value = 1.66
$ws.Range("U66").Value = 1.51
$ws.Range("W66").Value = 0.99
$ws.Range("Y66").Value = 0.85
$ws.Range("Z66").Value = 0.37
$ws.Range("W68").Value = -132.91
$ws.Range("X68").Value = 4.21
$ws.Range("T69").Value = 49.17
$ws.Range("U69").Value = 49.26
$ws.Range("V69").Value = 35.41
$ws.Range("X69").Value = 35.68
$ws.Range("Y69").Value = 30.78
$ws.Range("Z69").Value = 15.33
$ws.Range("T71").Value = 2.11
$ws.Range("U71").Value = 2.02
$ws.Range("V71").Value = 0.78
$ws.Range("W71").Value = 1.23
$ws.Range("Y71").Value = 1.11
$ws.Range("Z71").Value = 0.4
$ws.Range("W73").Value = -133.14
$ws.Range("X73").Value = 4.21
$ws.Range("U74").Value = 47.24
$ws.Range("V74").Value = 34.63
$ws.Range("X74").Value = 34.2
$ws.Range("Y74").Value = 29.67
$ws.Range("Z74").Value = 14.93
$ws.Range("W78").Value = -131.92
$ws.Range("X78").Value = 4.21
$ws.Range("U79").Value = 47.24
$ws.Range("V79").Value = 34.63
$ws.Range("X79").Value = 34.2
$ws.Range("Y79").Value = 29.67
$ws.Range("Z79").Value = 14.93
$ws.Range("W83").Value = -131.92
$ws.Range("X83").Value = 4.21
$ws.Range("U84").Value = 43.38
$ws.Range("V84").Value = 31.86
$ws.Range("X84").Value = 32.26
$ws.Range("Y84").Value = 29.32
$ws.Range("Z84").Value = 14.81
$ws.Range("U86").Value = -3.86
$ws.Range("V86").Value = -2.77
$ws.Range("W86").Value = -2.17
$ws.Range("Y86").Value = -0.35
$ws.Range("Z86").Value = -0.12
$ws.Range("W88").Value = -129.75
$ws.Range("X88").Value = 4.21
$ws.Range("U89").Value = 46.92
$ws.Range("V89").Value = 34.8
$ws.Range("Y89").Value = 30.27
$ws.Range("Z89").Value = 15.33
$ws.Range("U91").Value = -0.33
$ws.Range("V91").Value = 0.17
$ws.Range("W91").Value = 0.23
$ws.Range("X91").Value = 0.8
$ws.Range("Y91").Value = 0.61
$ws.Range("Z91").Value = 0.4
$ws.Range("W93").Value = -132.14
$ws.Range("X93").Value = 4.21
